# Updated cryptos list (price + 1h volume change columns) as published by the
# scheduled GitHub Actions scraper run. Values that look like a single plain
# number (e.g. "351.77") are written with a leading apostrophe so Excel keeps
# them as text instead of silently converting them to a numeric cell - this
# matches the original sheet where every Price/Volume cell is stored as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.185.99'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').Value = '2.886.26'
$ws.Range('E3').Value = '  +4.04%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''351.77'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').Value = '''111.73'
$ws.Range('E6').Value = '  +3.28%  '
$ws.Range('E7').Value = '  +1.99%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = '''0.623'
$ws.Range('E9').Value = '  +0.54%  '
$ws.Range('D10').Value = '''40.19'
$ws.Range('E10').Value = '  +2.43%  '
$ws.Range('D11').Value = '''0.0861'
$ws.Range('E11').Value = '  +3.33%  '
$ws.Range('D12').Value = '''0.135'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').Value = '''20.09'
$ws.Range('E13').Value = '  +1.34%  '
$ws.Range('D14').Value = '''7.85'
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('D15').Value = '3.338.16'
$ws.Range('E15').Value = '  +3.85%  '
$ws.Range('D16').Value = '''0.995'
$ws.Range('E16').Value = '  +7.78%  '
$ws.Range('D17').Value = '2.874.96'
$ws.Range('E17').Value = '  +2.73%  '
$ws.Range('D18').Value = '52.163.58'
$ws.Range('E18').Value = '  +1.43%  '
$ws.Range('E19').Value = '  +8.81%  '
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('D21').Value = '''13.65'
$ws.Range('E21').Value = '  +2.80%  '
$ws.Range('D22').Value = '0.0₃0983'
$ws.Range('E22').Value = '  +1.91%  '
$ws.Range('D23').Value = '''71.09'
$ws.Range('E23').Value = '  +0.91%  '
$ws.Range('D24').Value = '''269.49'
$ws.Range('E24').Value = '  +1.14%  '
$ws.Range('D25').Value = '''2.78'
$ws.Range('E25').Value = '  +1.56%  '
$ws.Range('D26').Value = '''26.43'
$ws.Range('E26').Value = '  +2.16%  '
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('D28').Value = '''0.164'
$ws.Range('E28').Value = '  +0.55%  '
$ws.Range('E29').Value = '  +3.31%  '
$ws.Range('D30').Value = '''38.80'
$ws.Range('E30').Value = '  +5.17%  '
$ws.Range('E31').Value = '  +1.09%  '
$ws.Range('D32').Value = '''6.35'
$ws.Range('E32').Value = '  +2.93%  '
$ws.Range('D33').Value = '''53.13'
$ws.Range('E33').Value = '  +2.13%  '
$ws.Range('D34').Value = '''5.97'
$ws.Range('E34').Value = '  +7.71%  '
$ws.Range('D35').Value = '''0.0925'
$ws.Range('E35').Value = '  +10.57%  '
$ws.Range('D36').Value = '''0.0458'
$ws.Range('E36').Value = '  +3.32%  '
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('E38').Value = '  +7.39%  '
$ws.Range('D39').Value = '''18.67'
$ws.Range('E39').Value = '  +0.68%  '
$ws.Range('E40').Value = '  +4.32%  '
$ws.Range('D41').Value = '''2.62'
$ws.Range('E41').Value = '  +5.40%  '
$ws.Range('D43').Value = '''22.70'
$ws.Range('E43').Value = '  +3.16%  '
$ws.Range('D44').Value = '''121.61'
$ws.Range('E44').Value = '  +0.96%  '
$ws.Range('E45').Value = '  +1.38%  '
$ws.Range('E46').Value = '  +8.76%  '
$ws.Range('D47').Value = '2.188.05'
$ws.Range('E47').Value = '  +2.31%  '
$ws.Range('E48').Value = '  +5.93%  '
$ws.Range('D49').Value = '''0.257'
$ws.Range('E49').Value = '  +15.23%  '
$ws.Range('D50').Value = '''0.957'
$ws.Range('E50').Value = '  +7.13%  '
$ws.Range('D51').Value = '''0.0323'
$ws.Range('E51').Value = '  +11.40%  '
